$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 becomes a date (2026-02-01, Excel serial 46054) formatted as "d-mmm"
$ws.Range("C2").NumberFormat = "d-mmm"
$ws.Range("C2").Value = 46054

# Append the new ranking rows (420-438)
$ws.Range("A420").Value = "הגר אגמון"
$ws.Range("B420").Value = 1

$ws.Range("A421").Value = "תאיו ורד"
$ws.Range("B421").Value = 1

$ws.Range("A422").Value = "מעיין סטרוזר"
$ws.Range("B422").Value = 1

$ws.Range("A423").Value = "תומר ששון"
$ws.Range("B423").Value = 1

$ws.Range("A424").Value = "שלו דיין"
$ws.Range("B424").Value = 1

$ws.Range("A425").Value = "מעיין סטרוזר"
$ws.Range("B425").Value = 6

$ws.Range("A426").Value = "שלו דיין"
$ws.Range("B426").Value = 6

$ws.Range("A427").Value = "לידור אלשטיין"
$ws.Range("B427").Value = 1

$ws.Range("A428").Value = "רומי הרשקוביץ"
$ws.Range("B428").Value = 1

$ws.Range("A429").Value = "איתי הראל"
$ws.Range("B429").Value = 1

$ws.Range("A430").Value = "ליהי בראל"
$ws.Range("B430").Value = 1

$ws.Range("A431").Value = "קרן רינת פביאן"
$ws.Range("B431").Value = 1

$ws.Range("A432").Value = "הילה שולויס"
$ws.Range("B432").Value = 1

$ws.Range("A433").Value = "אורי שטרנברג"
$ws.Range("B433").Value = 1

$ws.Range("A434").Value = "ירון גלפנד"
$ws.Range("B434").Value = 1

$ws.Range("A435").Value = "ליאם דיין"
$ws.Range("B435").Value = 1

$ws.Range("A436").Value = "איתי בסטקר"
$ws.Range("B436").Value = 1

$ws.Range("A437").Value = "אורי שטרנברג"
$ws.Range("B437").Value = 6

$ws.Range("A438").Value = "ירון גלפנד"
$ws.Range("B438").Value = 6

# Match the author's final selection state
$ws.Range("C423").Select()
